$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elements")
$x = "0" + ""
$ws.Range("Z1").Value = $x
